$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 78622
$ws.Range("C2").Value = 5444.512500000001
$ws.Range("D2").Value = 73177.4875

$ws.Range("B3").Value = 75199
$ws.Range("C3").Value = 5387.3435
$ws.Range("D3").Value = 69811.6565

$ws.Range("B4").Value = 69820
$ws.Range("C4").Value = 5357.2225
$ws.Range("D4").Value = 64462.7775

$ws.Range("B5").Value = 69709
$ws.Range("C5").Value = 5289.154500000001
$ws.Range("D5").Value = 64419.8455

$ws.Range("B6").Value = 69148
$ws.Range("C6").Value = 5361.8285
$ws.Range("D6").Value = 63786.1715

$ws.Range("B7").Value = 71026
$ws.Range("C7").Value = 5442.6015
$ws.Range("D7").Value = 65583.3985

$ws.Range("B8").Value = 69369
$ws.Range("C8").Value = 5816.530999999999
$ws.Range("D8").Value = 63552.469

$ws.Range("B9").Value = 83438
$ws.Range("C9").Value = 7034.6115
$ws.Range("D9").Value = 76403.3885

$ws.Range("B10").Value = 67948
$ws.Range("C10").Value = 8596.6895
$ws.Range("D10").Value = 59351.3105

$ws.Range("B11").ClearContents()
$ws.Range("C11").Value = 13942.873
$ws.Range("D11").Value = 99192.12700000001

$ws.Range("B12").ClearContents()
$ws.Range("C12").Value = 15769.3095
$ws.Range("D12").Value = 99767.6905

$ws.Range("C13").Value = 15209.229
$ws.Range("D13").Value = 100144.771

$ws.Range("C14").Value = 15159.795
$ws.Range("D14").Value = 101117.205

$ws.Range("C15").Value = 15449.4445
$ws.Range("D15").Value = 108724.5555

$ws.Range("C16").Value = 15588.433
$ws.Range("D16").Value = 103763.567

$ws.Range("C17").Value = 16020.6585
$ws.Range("D17").Value = 101211.3415

$ws.Range("C18").Value = 16262.3825
$ws.Range("D18").Value = 93839.6175

$ws.Range("C19").Value = 15862.609
$ws.Range("D19").Value = 88583.391

$ws.Range("C20").Value = 14887.418
$ws.Range("D20").Value = 89024.582

$ws.Range("C21").Value = 13341.0375
$ws.Range("D21").Value = 87063.9625

$ws.Range("C22").Value = 11837.672
$ws.Range("D22").Value = 85910.328

$ws.Range("C23").Value = 9490.991999999998
$ws.Range("D23").Value = 84684.008

$ws.Range("C24").Value = 6891.794
$ws.Range("D24").Value = 81619.206

$ws.Range("C25").Value = 5770.9715
$ws.Range("D25").Value = 52348.0285
